# Jake_Progress_Report.xlsx update
# - Hours column (B) switched from text descriptions ("45 minutes", "1 hour", "2 hours")
#   to numeric hour values for the existing rows.
# - New progress-report rows (6-9) added with dates, hours, and activity descriptions.
# - Active cell / selection moved to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5: Hours column becomes numeric ---
$ws.Range("B2").Value = 0.75
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 0.75

# --- Fill in new rows 6-9 ---
# Copy the date formatting (numFmt + border) from an existing populated date cell (A3)
# so the newly-entered dates pick up the same display style instead of a brand new one.
$ws.Range("A3").Copy()
$ws.Range("A6:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A6").Value = 42793
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "Enhanced the category tree"

$ws.Range("A7").Value = 42795
$ws.Range("B7").Value = 1.5
$ws.Range("C7").Value = "Gathering Data for registered user's table"

$ws.Range("A8").Value = 42796
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Designing conceptual prototype for website"

$ws.Range("A9").Value = 42796
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "Writing SQL to populate Sales_items, Suppliers, Users, Transactions tables"

# --- Update selection to C9 ---
$ws.Range("C9").Select() | Out-Null
